$d = $word.ActiveDocument

# The paragraph that contains the sentence we need to split/expand.
# "Describe the trait, the index, and the index equation. Include the
# reference." is entirely contained in a single run; we need to turn it
# into five runs with new wording while leaving every other run in the
# document untouched.
$oldSentence = "Describe the trait, the index, and the index equation. Include the reference."

$hit = $d.Content
$hit.Find.ClearFormatting()
$hit.Find.Forward = $true
$found = $hit.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence."
}

# Rebuild a plain Range from the match bounds: InsertXML needs a Range
# object created directly via Document.Range to reliably replace content
# in place (re-using the Find hit range directly was unreliable). This
# range spans exactly the old sentence; InsertXML replaces its content in
# one shot with the OOXML we supply, so adjacent runs elsewhere in the
# paragraph are left completely alone (no implicit re-merge of
# neighbouring runs that happen to share identical formatting).
$rng = $d.Range($hit.Start, $hit.End)
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

$newRuns = (
    '<w:r w:rsidR="00212854">' + $rPr + '<w:t xml:space="preserve">Describe the trait, the index, </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">the </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>index equation</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>, and the expected parametric space (observed range)</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>. Include the reference.</w:t></w:r>'
)

$packageXml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $newRuns + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
)

$rng.InsertXML($packageXml)
